$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap order of "Malaui" and "Mozambique" rows (row 111 was Malaui / row 112 was
# Mozambique; after the update row 111 is Mozambique with refreshed figures and row 112
# becomes Malaui carrying the previous Malaui figures) ---
$ws.Range("A111").Value = "Mozambique"
$ws.Range("A112").Value = "Malaui"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 19:20"

# --- Update country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6758709
$ws.Range("C4").Value = 9420
$ws.Range("D4").Value = 4040498
$ws.Range("E4").Value = 2518801
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 410
$ws.Range("H4").Value = 199410

# Row 5: India
$ws.Range("B5").Value = 5009290
$ws.Range("C5").Value = 82376
$ws.Range("D5").Value = 3933455
$ws.Range("E5").Value = 993790
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1237
$ws.Range("H5").Value = 82045

# Row 6: Brasil
$ws.Range("B6").Value = 4356690
$ws.Range("C6").Value = 7146
$ws.Range("D6").Value = 3613184
$ws.Range("E6").Value = 611209
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 180
$ws.Range("H6").Value = 132297

# Row 12: España
$ws.Range("B12").Value = 603167
$ws.Range("C12").Value = 9437
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 156
$ws.Range("H12").Value = 30004

# Row 21: Irak
$ws.Range("B21").Value = 298702
$ws.Range("C21").Value = 4224
$ws.Range("D21").Value = 233346
$ws.Range("E21").Value = 57190
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 8166

# Row 22: Turquia
$ws.Range("B22").Value = 294620
$ws.Range("C22").Value = 1742
$ws.Range("D22").Value = 261260
$ws.Range("E22").Value = 26174
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 67
$ws.Range("H22").Value = 7186

# Row 25: Alemania
$ws.Range("B25").Value = 264169
$ws.Range("C25").Value = 948
$ws.Range("D25").Value = 237550
$ws.Range("E25").Value = 17177
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9442

# Row 29: Canada
$ws.Range("B29").Value = 138555
$ws.Range("C29").Value = 545
$ws.Range("D29").Value = 121224
$ws.Range("E29").Value = 8152
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 9179

# Row 32: Ecuador
$ws.Range("B32").Value = 119553
$ws.Range("C32").Value = 642
$ws.Range("D32").Value = 97063
$ws.Range("E32").Value = 11527
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 41
$ws.Range("H32").Value = 10963

# Row 59: Argelia
$ws.Range("B59").Value = 48734
$ws.Range("C59").Value = 238
$ws.Range("D59").Value = 34385
$ws.Range("E59").Value = 12717
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 12
$ws.Range("H59").Value = 1632

# Row 68: Chequia
$ws.Range("B68").Value = 38187
$ws.Range("C68").Value = 965
$ws.Range("D68").Value = 22526
$ws.Range("E68").Value = 15185
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 11
$ws.Range("H68").Value = 476

# Row 77: Libano
$ws.Range("B77").Value = 25449
$ws.Range("C77").Value = 592
$ws.Range("D77").Value = 8765
$ws.Range("E77").Value = 16432
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 6
$ws.Range("H77").Value = 252

# Row 100: Maldivas
$ws.Range("B100").Value = 9328
$ws.Range("C100").Value = 85
$ws.Range("D100").Value = 7729
$ws.Range("E100").Value = 1566
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 33

# Row 111: Mozambique (new figures)
$ws.Range("B111").Value = 5713
$ws.Range("C111").Value = 231
$ws.Range("D111").Value = 3181
$ws.Range("E111").Value = 2495
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 37

# Row 112: Malaui (previous Malaui figures)
$ws.Range("B112").Value = 5697
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 3742
$ws.Range("E112").Value = 1777
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 178

# Row 170: San Marino
$ws.Range("B170").Value = 723
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 663
$ws.Range("E170").Value = 18
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 42

# Row 179: Islas Feroe
$ws.Range("B179").Value = 428
$ws.Range("C179").Value = 5
$ws.Range("D179").Value = 412
$ws.Range("E179").Value = 16
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0
